# Rename worksheets to new summ<ID> identifiers as part of rerunning the
# LU d2c FeatEng pipeline for FR cities with new spatial units / dist models.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ12186826",
    "summ12667636",
    "summ13167217",
    "summ13649434",
    "summ14208109",
    "summ14700449",
    "summ15200074",
    "summ15866536",
    "summ16802683"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
